# Daily update at 8 AM UTC
# Adds the next day's row (day 45965) to the Wins Over Time log and
# moves the "latest row" date-only formatting down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-last row (15) loses its date-only ("YYYY-MM-DD") look and
# goes back to the regular date+time format used by every other data row.
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 16.
$ws.Range("A16").Value = 45965
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = 41
$ws.Range("D16").Value = 38

# The new last row takes on the date-only formatting.
$ws.Range("A16").NumberFormat = "YYYY-MM-DD"
